$d = $word.ActiveDocument

# Insert the new paragraph after paragraph 27 ("...scorrevolezza.")
$p27 = $d.Paragraphs(27).Range
$p27.InsertParagraphAfter()
$newPara = $d.Paragraphs(28)
$newPara.Range.Text = "Per quanto riguarda l’audio, si sono riscontrati degli errori nell’utilizzare la libreria di default di Processing. Dopo una ricerca sul forum di Processing, si è optato per l’utilizzo della libreria Minim, con cui si è ottenuto lo stesso risultato."

# Move the _GoBack bookmark to the end of the new paragraph (empty range right before the paragraph mark)
$newParaRange = $d.Paragraphs(28).Range
$bmPos = $newParaRange.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "Done"
